$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Range("E10").Value = "waitType::visible" + [char]10 + "export::abc::WebElement::getText"
$ws.Range("E10").Style = $ws.Range("E3").Style
$ws.Range("F10").Value = ""
Write-Host "done"
